$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = [double]"12.265976"
$ws.Cells.Item(2, 8).Value = [double]"36.797928"
$ws.Cells.Item(2, 9).Value = [double]"0.004000867643088759"
$ws.Cells.Item(2, 10).Value = [double]"0.004000867643088758"
$ws.Cells.Item(2, 13).Value = [double]"2.685464"
$ws.Cells.Item(2, 14).Value = [double]"8.056392000000001"
$ws.Cells.Item(2, 15).Value = [double]"0.06676031826184478"
$ws.Cells.Item(2, 16).Value = [double]"0.06676031826184478"
$ws.Cells.Item(2, 17).Value = [double]"32.939836972864"
$ws.Cells.Item(2, 18).Value = [double]"296.458532755776"
$ws.Cells.Item(2, 19).Value = [double]"0.0002670991971761224"
$ws.Cells.Item(2, 20).Value = [double]"0.0002670991971761224"

$ws.Cells.Item(3, 7).Value = [double]"12.265976"
$ws.Cells.Item(3, 8).Value = [double]"36.797928"
$ws.Cells.Item(3, 9).Value = [double]"0.004000867643088759"
$ws.Cells.Item(3, 10).Value = [double]"0.004000867643088758"
$ws.Cells.Item(3, 15).Value = [double]"0.02342101692711854"
$ws.Cells.Item(3, 16).Value = [double]"0.02342101692711854"
$ws.Cells.Item(3, 17).Value = [double]"11.556033575096"
$ws.Cells.Item(3, 18).Value = [double]"104.004302175864"
$ws.Cells.Item(3, 19).Value = [double]"9.370438879194269E-05"
$ws.Cells.Item(3, 20).Value = [double]"9.370438879194267E-05"

$ws.Cells.Item(4, 7).Value = [double]"12.265976"
$ws.Cells.Item(4, 8).Value = [double]"36.797928"
$ws.Cells.Item(4, 9).Value = [double]"0.004000867643088759"
$ws.Cells.Item(4, 10).Value = [double]"0.004000867643088758"
$ws.Cells.Item(4, 13).Value = [double]"35.399925"
$ws.Cells.Item(4, 14).Value = [double]"106.199775"
$ws.Cells.Item(4, 15).Value = [double]"0.8800379597140142"
$ws.Cells.Item(4, 16).Value = [double]"0.8800379597140142"
$ws.Cells.Item(4, 17).Value = [double]"434.2146304518"
$ws.Cells.Item(4, 18).Value = [double]"3907.9316740662"
$ws.Cells.Item(4, 19).Value = [double]"0.003520915397709649"
$ws.Cells.Item(4, 20).Value = [double]"0.003520915397709648"

$ws.Cells.Item(5, 7).Value = [double]"12.265976"
$ws.Cells.Item(5, 8).Value = [double]"36.797928"
$ws.Cells.Item(5, 9).Value = [double]"0.004000867643088759"
$ws.Cells.Item(5, 10).Value = [double]"0.004000867643088758"
$ws.Cells.Item(5, 13).Value = [double]"1.197942333333333"
$ws.Cells.Item(5, 14).Value = [double]"3.593827"
$ws.Cells.Item(5, 15).Value = [double]"0.02978070509702244"
$ws.Cells.Item(5, 16).Value = [double]"0.02978070509702244"
$ws.Cells.Item(5, 17).Value = [double]"14.69393191005067"
$ws.Cells.Item(5, 18).Value = [double]"132.245387190456"
$ws.Cells.Item(5, 19).Value = [double]"0.0001191486594110456"
$ws.Cells.Item(5, 20).Value = [double]"0.0001191486594110456"

$ws.Cells.Item(6, 9).Value = [double]"0.9924545876219728"
$ws.Cells.Item(6, 10).Value = [double]"0.9924545876219727"
$ws.Cells.Item(6, 13).Value = [double]"2.685464"
$ws.Cells.Item(6, 14).Value = [double]"8.056392000000001"
$ws.Cells.Item(6, 15).Value = [double]"0.06676031826184478"
$ws.Cells.Item(6, 16).Value = [double]"0.06676031826184478"
$ws.Cells.Item(6, 17).Value = [double]"8171.050690894725"
$ws.Cells.Item(6, 18).Value = [double]"73539.45621805252"
$ws.Cells.Item(6, 19).Value = [double]"0.06625658413007082"
$ws.Cells.Item(6, 20).Value = [double]"0.06625658413007081"

$ws.Cells.Item(7, 9).Value = [double]"0.9924545876219728"
$ws.Cells.Item(7, 10).Value = [double]"0.9924545876219727"
$ws.Cells.Item(7, 15).Value = [double]"0.02342101692711854"
$ws.Cells.Item(7, 16).Value = [double]"0.02342101692711854"
$ws.Cells.Item(7, 19).Value = [double]"0.02324429569609068"
$ws.Cells.Item(7, 20).Value = [double]"0.02324429569609067"

$ws.Cells.Item(8, 9).Value = [double]"0.9924545876219728"
$ws.Cells.Item(8, 10).Value = [double]"0.9924545876219727"
$ws.Cells.Item(8, 13).Value = [double]"35.399925"
$ws.Cells.Item(8, 14).Value = [double]"106.199775"
$ws.Cells.Item(8, 15).Value = [double]"0.8800379597140142"
$ws.Cells.Item(8, 16).Value = [double]"0.8800379597140142"
$ws.Cells.Item(8, 17).Value = [double]"107711.2117789967"
$ws.Cells.Item(8, 18).Value = [double]"969400.9060109699"
$ws.Cells.Item(8, 19).Value = [double]"0.8733977103996543"
$ws.Cells.Item(8, 20).Value = [double]"0.8733977103996542"

$ws.Cells.Item(9, 9).Value = [double]"0.9924545876219728"
$ws.Cells.Item(9, 10).Value = [double]"0.9924545876219727"
$ws.Cells.Item(9, 13).Value = [double]"1.197942333333333"
$ws.Cells.Item(9, 14).Value = [double]"3.593827"
$ws.Cells.Item(9, 15).Value = [double]"0.02978070509702244"
$ws.Cells.Item(9, 16).Value = [double]"0.02978070509702244"
$ws.Cells.Item(9, 17).Value = [double]"3644.974399372091"
$ws.Cells.Item(9, 18).Value = [double]"32804.76959434881"
$ws.Cells.Item(9, 19).Value = [double]"0.02955599739615699"
$ws.Cells.Item(9, 20).Value = [double]"0.02955599739615699"

$ws.Cells.Item(10, 7).Value = [double]"8.377189333333332"
$ws.Cells.Item(10, 8).Value = [double]"25.131568"
$ws.Cells.Item(10, 9).Value = [double]"0.002732438555542716"
$ws.Cells.Item(10, 10).Value = [double]"0.002732438555542716"
$ws.Cells.Item(10, 13).Value = [double]"2.685464"
$ws.Cells.Item(10, 14).Value = [double]"8.056392000000001"
$ws.Cells.Item(10, 15).Value = [double]"0.06676031826184478"
$ws.Cells.Item(10, 16).Value = [double]"0.06676031826184478"
$ws.Cells.Item(10, 17).Value = [double]"22.49664037585066"
$ws.Cells.Item(10, 18).Value = [double]"202.469763382656"
$ws.Cells.Item(10, 19).Value = [double]"0.0001824184675989671"
$ws.Cells.Item(10, 20).Value = [double]"0.0001824184675989671"

$ws.Cells.Item(11, 7).Value = [double]"8.377189333333332"
$ws.Cells.Item(11, 8).Value = [double]"25.131568"
$ws.Cells.Item(11, 9).Value = [double]"0.002732438555542716"
$ws.Cells.Item(11, 10).Value = [double]"0.002732438555542716"
$ws.Cells.Item(11, 15).Value = [double]"0.02342101692711854"
$ws.Cells.Item(11, 16).Value = [double]"0.02342101692711854"
$ws.Cells.Item(11, 17).Value = [double]"7.892325991909333"
$ws.Cells.Item(11, 18).Value = [double]"71.03093392718399"
$ws.Cells.Item(11, 19).Value = [double]"6.399648966167728E-05"
$ws.Cells.Item(11, 20).Value = [double]"6.399648966167728E-05"

$ws.Cells.Item(12, 7).Value = [double]"8.377189333333332"
$ws.Cells.Item(12, 8).Value = [double]"25.131568"
$ws.Cells.Item(12, 9).Value = [double]"0.002732438555542716"
$ws.Cells.Item(12, 10).Value = [double]"0.002732438555542716"
$ws.Cells.Item(12, 13).Value = [double]"35.399925"
$ws.Cells.Item(12, 14).Value = [double]"106.199775"
$ws.Cells.Item(12, 15).Value = [double]"0.8800379597140142"
$ws.Cells.Item(12, 16).Value = [double]"0.8800379597140142"
$ws.Cells.Item(12, 17).Value = [double]"296.5518741108"
$ws.Cells.Item(12, 18).Value = [double]"2668.9668669972"
$ws.Cells.Item(12, 19).Value = [double]"0.00240464965146372"
$ws.Cells.Item(12, 20).Value = [double]"0.00240464965146372"

$ws.Cells.Item(13, 7).Value = [double]"8.377189333333332"
$ws.Cells.Item(13, 8).Value = [double]"25.131568"
$ws.Cells.Item(13, 9).Value = [double]"0.002732438555542716"
$ws.Cells.Item(13, 10).Value = [double]"0.002732438555542716"
$ws.Cells.Item(13, 13).Value = [double]"1.197942333333333"
$ws.Cells.Item(13, 14).Value = [double]"3.593827"
$ws.Cells.Item(13, 15).Value = [double]"0.02978070509702244"
$ws.Cells.Item(13, 16).Value = [double]"0.02978070509702244"
$ws.Cells.Item(13, 17).Value = [double]"10.03538973674844"
$ws.Cells.Item(13, 18).Value = [double]"90.31850763073599"
$ws.Cells.Item(13, 19).Value = [double]"8.137394681835159E-05"
$ws.Cells.Item(13, 20).Value = [double]"8.137394681835159E-05"

$ws.Cells.Item(14, 7).Value = [double]"2.489778666666667"
$ws.Cells.Item(14, 8).Value = [double]"7.469336"
$ws.Cells.Item(14, 9).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(14, 10).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(14, 13).Value = [double]"2.685464"
$ws.Cells.Item(14, 14).Value = [double]"8.056392000000001"
$ws.Cells.Item(14, 15).Value = [double]"0.06676031826184478"
$ws.Cells.Item(14, 16).Value = [double]"0.06676031826184478"
$ws.Cells.Item(14, 17).Value = [double]"6.686210977301333"
$ws.Cells.Item(14, 18).Value = [double]"60.17589879571201"
$ws.Cells.Item(14, 19).Value = [double]"5.421646699886768E-05"
$ws.Cells.Item(14, 20).Value = [double]"5.421646699886768E-05"

$ws.Cells.Item(15, 7).Value = [double]"2.489778666666667"
$ws.Cells.Item(15, 8).Value = [double]"7.469336"
$ws.Cells.Item(15, 9).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(15, 10).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(15, 15).Value = [double]"0.02342101692711854"
$ws.Cells.Item(15, 16).Value = [double]"0.02342101692711854"
$ws.Cells.Item(15, 17).Value = [double]"2.345672767218667"
$ws.Cells.Item(15, 18).Value = [double]"21.111054904968"
$ws.Cells.Item(15, 19).Value = [double]"1.902035257424424E-05"
$ws.Cells.Item(15, 20).Value = [double]"1.902035257424424E-05"

$ws.Cells.Item(16, 7).Value = [double]"2.489778666666667"
$ws.Cells.Item(16, 8).Value = [double]"7.469336"
$ws.Cells.Item(16, 9).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(16, 10).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(16, 13).Value = [double]"35.399925"
$ws.Cells.Item(16, 14).Value = [double]"106.199775"
$ws.Cells.Item(16, 15).Value = [double]"0.8800379597140142"
$ws.Cells.Item(16, 16).Value = [double]"0.8800379597140142"
$ws.Cells.Item(16, 17).Value = [double]"88.1379780666"
$ws.Cells.Item(16, 18).Value = [double]"793.2418025994"
$ws.Cells.Item(16, 19).Value = [double]"0.0007146842651865342"
$ws.Cells.Item(16, 20).Value = [double]"0.0007146842651865342"

$ws.Cells.Item(17, 7).Value = [double]"2.489778666666667"
$ws.Cells.Item(17, 8).Value = [double]"7.469336"
$ws.Cells.Item(17, 9).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(17, 10).Value = [double]"0.0008121061793956991"
$ws.Cells.Item(17, 13).Value = [double]"1.197942333333333"
$ws.Cells.Item(17, 14).Value = [double]"3.593827"
$ws.Cells.Item(17, 15).Value = [double]"0.02978070509702244"
$ws.Cells.Item(17, 16).Value = [double]"0.02978070509702244"
$ws.Cells.Item(17, 17).Value = [double]"2.982611265430223"
$ws.Cells.Item(17, 18).Value = [double]"26.843501388872"
$ws.Cells.Item(17, 19).Value = [double]"2.418509463605292E-05"
$ws.Cells.Item(17, 20).Value = [double]"2.418509463605292E-05"
